$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 21, shifting existing rows 21-49 down to 23-51
$ws.Rows.Item(21).Resize(2).Insert()

# Fill in new row 21 (copy constant columns from the row above, then set varying values)
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).Value = 44495
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112021
$ws.Cells.Item(21, 7).Value = "Ají"
$ws.Cells.Item(21, 8).Value = "Americana (o)"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 21
$ws.Cells.Item(21, 11).Value = 47000
$ws.Cells.Item(21, 12).Value = 48000
$ws.Cells.Item(21, 13).Value = 47571
$ws.Cells.Item(21, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(21, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(21, 16).Value = 1903
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"

# Fill in new row 22
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44495
$ws.Cells.Item(22, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 100112021
$ws.Cells.Item(22, 7).Value = "Ají"
$ws.Cells.Item(22, 8).Value = "Inferno"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 15
$ws.Cells.Item(22, 11).Value = 21000
$ws.Cells.Item(22, 12).Value = 22000
$ws.Cells.Item(22, 13).Value = 21333
$ws.Cells.Item(22, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 1778
$ws.Cells.Item(22, 17).Value = 12
$ws.Cells.Item(22, 18).Value = "Hortaliza"

Write-Host "Done"
